# This sheet is a weekly price log. Two new weekly records need to be added
# for the Terminal Hortofrutícola Agro Chillán - Frutilla series, inserted
# right before the existing row 125 (shifting all rows from 125 on down by two).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 125; this pushes old rows 125-137 down to 127-139
# and keeps their contents and formatting untouched.
$ws.Range("A125:A126").EntireRow.Insert()

# New row 125: Primera quality, 160 volume, $/caja 7 kilos, Provincia de Diguillín
$ws.Cells.Item(125, 1).Value = 7
$ws.Cells.Item(125, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(125, 3).Value = "Ñuble"
$ws.Cells.Item(125, 4).Value = 44491
$ws.Cells.Item(125, 5).Value = 16
$ws.Cells.Item(125, 6).Value = "Fruta"
$ws.Cells.Item(125, 7).Value = 100101
$ws.Cells.Item(125, 8).Value = "Berries"
$ws.Cells.Item(125, 9).Value = 100112025
$ws.Cells.Item(125, 10).Value = "Frutilla"
$ws.Cells.Item(125, 11).Value = "Sin especificar"
$ws.Cells.Item(125, 12).Value = "Primera"
$ws.Cells.Item(125, 13).Value = 160
$ws.Cells.Item(125, 14).Value = 7500
$ws.Cells.Item(125, 15).Value = 8000
$ws.Cells.Item(125, 16).Value = 7750
$ws.Cells.Item(125, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(125, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(125, 19).Value = 1107
$ws.Cells.Item(125, 20).Value = 7

# New row 126: Segunda quality, 160 volume, $/caja 7 kilos, Provincia de Diguillín
$ws.Cells.Item(126, 1).Value = 7
$ws.Cells.Item(126, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(126, 3).Value = "Ñuble"
$ws.Cells.Item(126, 4).Value = 44491
$ws.Cells.Item(126, 5).Value = 16
$ws.Cells.Item(126, 6).Value = "Fruta"
$ws.Cells.Item(126, 7).Value = 100101
$ws.Cells.Item(126, 8).Value = "Berries"
$ws.Cells.Item(126, 9).Value = 100112025
$ws.Cells.Item(126, 10).Value = "Frutilla"
$ws.Cells.Item(126, 11).Value = "Sin especificar"
$ws.Cells.Item(126, 12).Value = "Segunda"
$ws.Cells.Item(126, 13).Value = 160
$ws.Cells.Item(126, 14).Value = 6000
$ws.Cells.Item(126, 15).Value = 6500
$ws.Cells.Item(126, 16).Value = 6250
$ws.Cells.Item(126, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(126, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(126, 19).Value = 893
$ws.Cells.Item(126, 20).Value = 7
